# Revision restore: Rule "R30" From-value (Rules!C10) reverts 18 -> 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Cells.Item(10, 3).Value = 1
